$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E ("boson"), shifting pt_min's
# neighbours (boson..syst2_u) one column to the right, and add the new
# "pt_max" column (value 50 for every data row).
$ws.Range("E1").EntireColumn.Insert()

$ws.Range("E1").Value = "pt_max"
$ws.Range("E2:E12").Value = 50

$ws.Range("E2:E12").Select()
